# Clean services data and export as .csv
# - Normalize a handful of inconsistently-formatted phone numbers in column D
#   to match the "(XXX) XXX-XXXX" style used throughout the rest of the sheet.
# - Apply a Text ("@") number format to the Organization/Address/Phone columns
#   (B:D) for the header row and all data rows, so phone numbers such as area
#   codes/leading content are preserved as text rather than being reinterpreted.
# - Restore the view (zoom/selection) to match how the sheet was left after
#   the cleanup pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up inconsistent phone number formatting -----------------------
$ws.Range("D9").Value  = "(519) 434-2183"
$ws.Range("D18").Value = "(416) 554-1286"
$ws.Range("D48").Value = "(905) 876-2473"
$ws.Range("D50").Value = "(519) 822-2273"
$ws.Range("D51").Value = "(519) 794-9981"

# --- Apply Text number format to Organization/Address/Phone columns ----
$ws.Range("B1:D1").NumberFormat = "@"
$ws.Range("B2:D53").NumberFormat = "@"

# --- Restore view state --------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Zoom = 138
[void]$ws.Range("D3").Select()
